$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.580.02'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '2.294.93'
$ws.Range("E3").Value = '  +0.64%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = "'311.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.77%  '
$ws.Range("D6").Value = "'104.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.61%  '
$ws.Range("D7").Value = "'0.623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.41%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("D10").Value = "'39.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.72%  '
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("D12").Value = "'8.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.13%  '
$ws.Range("E13").Value = '  +1.90%  '
$ws.Range("E14").Value = '  +3.71%  '
$ws.Range("D15").Value = "'15.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("D16").Value = '2.643.00'
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").Value = '2.313.80'
$ws.Range("E17").Value = '  +1.17%  '
$ws.Range("D18").Value = '42.726.03'
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = "'13.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.42%  '
$ws.Range("D22").Value = "'73.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("D23").Value = "'3.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.99%  '
$ws.Range("D24").Value = "'263.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("D25").Value = "'2.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.51%  '
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("D27").Value = "'10.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").Value = "'7.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +17.42%  '
$ws.Range("D29").Value = "'2.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").Value = "'22.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("D31").Value = "'35.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.68%  '
$ws.Range("D32").Value = "'165.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("E33").Value = '  -0.30%  '
$ws.Range("D34").Value = "'0.129"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.03%  '
$ws.Range("D35").Value = "'2.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.44%  '
$ws.Range("E36").Value = '  -1.70%  '
$ws.Range("D37").Value = "'4.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.03%  '
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("E39").Value = '  +3.16%  '
$ws.Range("D40").Value = "'2.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.48%  '
$ws.Range("D41").Value = "'1.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.11%  '
$ws.Range("D42").Value = "'99.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.20%  '
$ws.Range("D43").Value = "'69.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.84%  '
$ws.Range("E44").Value = '  +1.74%  '
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("D46").Value = '1.742.55'
$ws.Range("E46").Value = '  +9.02%  '
$ws.Range("D47").Value = "'12.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.71%  '
$ws.Range("D48").Value = "'79.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.59%  '
$ws.Range("D49").Value = "'110.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.05%  '
$ws.Range("D50").Value = "'5.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("E51").Value = '  -2.65%  '
